# Updates cryptos list figures (prices, volume %, and a 3-row reorder)
# to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.029.44"
$ws.Range("E2").Value = "  +1.21%  "

$ws.Range("D3").Value = "1.643.46"
$ws.Range("E3").Value = "  +2.54%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.67%  "

$ws.Range("E6").Value = "  +1.02%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.40%  "

$ws.Range("E9").Value = "  +3.33%  "

$ws.Range("E10").Value = "  +1.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0917"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "1.877.73"
$ws.Range("E12").Value = "  +2.51%  "

$ws.Range("D13").Value = "1.644.18"
$ws.Range("E13").Value = "  +2.50%  "

$ws.Range("E14").Value = "  +5.19%  "

$ws.Range("E15").Value = "  +21.37%  "

$ws.Range("E16").Value = "  +4.73%  "

$ws.Range("D17").Value = "30.073.93"
$ws.Range("E17").Value = "  +1.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.75%  "

$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("E27").Value = "  +2.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.94%  "

$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0494"
$ws.Range("D30").Style = "Normal"

$ws.Range("E31").Value = "  +6.22%  "

$ws.Range("E32").Value = "  +5.90%  "

$ws.Range("E33").Value = "  +0.37%  "

$ws.Range("D34").Value = "1.438.26"
$ws.Range("E34").Value = "  +0.87%  "

$ws.Range("E35").Value = "  +7.86%  "

$ws.Range("E36").Value = "  +0.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.97%  "

$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0172"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.58%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "77.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.29%  "

$ws.Range("E41").Value = "  +2.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.45%  "

$ws.Range("E43").Value = "  +3.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "54.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.89%  "

$ws.Range("E45").Value = "  -0.58%  "

$ws.Range("E46").Value = "  +6.97%  "

$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("D49").Value = "1.784.58"
$ws.Range("E49").Value = "  +2.41%  "

$ws.Range("E50").Value = "  +9.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.99%  "
